# autoBW workbook update:
#  - "Create Activities": switch to the new activity-import column layout
#    (database/activity_name/reference_product/reference_product_amount/
#     reference_product_unit/std_dev/location/activity_version/activity_code)
#    and populate it with the grid-mix / fancy-chemical / 2-Methyl pentane
#    activities.
#  - "Add Exchanges": switch to the new exchange-import column layout
#    (database/activity_name/exchange_db/exchange/amount/unit/location)
#    and populate it with the exchanges wiring those activities together.
#  - "Delete Exchanges": unchanged data, selection only.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Create Activities"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Create Activities")

$ws1.Range("A1:J2").ClearContents()

# Header row
$ws1.Cells.Item(1,1).Value = "database"
$ws1.Cells.Item(1,2).Value = "activity_name"
$ws1.Cells.Item(1,3).Value = "reference_product"
$ws1.Cells.Item(1,4).Value = "reference_product_amount"
$ws1.Cells.Item(1,5).Value = "reference_product_unit"
$ws1.Cells.Item(1,6).Value = "std_dev"
$ws1.Cells.Item(1,7).Value = "location"
$ws1.Cells.Item(1,8).Value = "activity_version"
$ws1.Cells.Item(1,9).Value = "activity_code"

# Row 2: 2-Methyl pentane, from chemistry, at plant
$ws1.Cells.Item(2,1).Value = "newdb"
$ws1.Cells.Item(2,2).Value = "2-Methyl pentane,  from chemistry"
$ws1.Cells.Item(2,3).Value = "2-Methyl pentane, from chemistry, at plant"
$ws1.Cells.Item(2,4).Value = 1
$ws1.Cells.Item(2,5).Value = "kg"
$ws1.Cells.Item(2,6).Value = 2
$ws1.Cells.Item(2,7).Value = "US"
$ws1.Cells.Item(2,8).Value = 0.1

# Row 3: fancy chemical, from a technology, at plant
$ws1.Cells.Item(3,1).Value = "newdb"
$ws1.Cells.Item(3,2).Value = "fancy chemical, from a technology, at plant"
$ws1.Cells.Item(3,3).Value = "fancy chemical"
$ws1.Cells.Item(3,4).Value = 1
$ws1.Cells.Item(3,5).Value = "kg"
$ws1.Cells.Item(3,7).Value = "US"
$ws1.Cells.Item(3,8).Value = 0.1

# Row 4: electric grid mix
$ws1.Cells.Item(4,1).Value = "newdb"
$ws1.Cells.Item(4,2).Value = "electric grid mix"
$ws1.Cells.Item(4,3).Value = "electricity"
$ws1.Cells.Item(4,4).Value = 1
$ws1.Cells.Item(4,5).Value = "kWh"
$ws1.Cells.Item(4,7).Value = "US"
$ws1.Cells.Item(4,8).Value = 0.1

# Best-effort column widths (bestFit values Excel would compute on open)
$ws1.Columns.Item(1).ColumnWidth = 9
$ws1.Columns.Item(2).ColumnWidth = 32.43
$ws1.Columns.Item(3).ColumnWidth = 40
$ws1.Columns.Item(4).ColumnWidth = 26
$ws1.Columns.Item(5).ColumnWidth = 22.57
$ws1.Columns.Item(6).ColumnWidth = 8
$ws1.Columns.Item(7).ColumnWidth = 8.14
$ws1.Columns.Item(8).ColumnWidth = 15.14
$ws1.Columns.Item(9).ColumnWidth = 12.71
$ws1.Columns.Item(10).Delete()

# ---------------------------------------------------------------------------
# Sheet 2: "Add Exchanges"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Add Exchanges")

$ws2.Range("A1:P8").ClearContents()

# Header row (bold)
$ws2.Cells.Item(1,1).Value = "database"
$ws2.Cells.Item(1,2).Value = "activity_name"
$ws2.Cells.Item(1,3).Value = "exchange_db"
$ws2.Cells.Item(1,4).Value = "exchange"
$ws2.Cells.Item(1,5).Value = "amount"
$ws2.Cells.Item(1,6).Value = "unit"
$ws2.Cells.Item(1,7).Value = "location"
$ws2.Range("A1:M1").Font.Bold = $true

# Row 2: electric grid mix <- 2-Methyl pentane, from chemistry
$ws2.Cells.Item(2,1).Value = "newdb"
$ws2.Cells.Item(2,2).Value = "electric grid mix"
$ws2.Cells.Item(2,3).Value = "newdb"
$ws2.Cells.Item(2,4).Value = "2-Methyl pentane, from chemistry"
$ws2.Cells.Item(2,5).Value = 0.8
$ws2.Cells.Item(2,6).Value = "kg"
$ws2.Cells.Item(2,7).Value = "US"

# Row 3: electric grid mix <- electricity
$ws2.Cells.Item(3,1).Value = "newdb"
$ws2.Cells.Item(3,2).Value = "electric grid mix"
$ws2.Cells.Item(3,3).Value = "newdb"
$ws2.Cells.Item(3,4).Value = "electricity"
$ws2.Cells.Item(3,5).Value = 0.1
$ws2.Cells.Item(3,6).Value = "kWh"
$ws2.Cells.Item(3,7).Value = "US"

# Row 4: 2-Methyl pentane, from chemistry, at plant <- electricity
$ws2.Cells.Item(4,1).Value = "newdb"
$ws2.Cells.Item(4,2).Value = "2-Methyl pentane,  from chemistry, at plant"
$ws2.Cells.Item(4,3).Value = "newdb"
$ws2.Cells.Item(4,4).Value = "electricity"
$ws2.Cells.Item(4,5).Value = 4
$ws2.Cells.Item(4,6).Value = "kWh"
$ws2.Cells.Item(4,7).Value = "US"

# Row 5: 2-Methyl pentane, from chemistry, at plant <- fancy chemical
$ws2.Cells.Item(5,1).Value = "newdb"
$ws2.Cells.Item(5,2).Value = "2-Methyl pentane,  from chemistry, at plant"
$ws2.Cells.Item(5,3).Value = "newdb"
$ws2.Cells.Item(5,4).Value = "fancy chemical"
$ws2.Cells.Item(5,5).Value = 0.2
$ws2.Cells.Item(5,6).Value = "kg"
$ws2.Cells.Item(5,7).Value = "US"

# Row 6: 2-Methyl pentane, from chemistry, at plant <- 2-Methyl pentane, from chemistry
$ws2.Cells.Item(6,1).Value = "newdb"
$ws2.Cells.Item(6,2).Value = "2-Methyl pentane,  from chemistry, at plant"
$ws2.Cells.Item(6,3).Value = "newdb"
$ws2.Cells.Item(6,4).Value = "2-Methyl pentane, from chemistry"
$ws2.Cells.Item(6,5).Value = 0.01
$ws2.Cells.Item(6,6).Value = "kg"
$ws2.Cells.Item(6,7).Value = "US"

# Row 7: fancy chemical, from a technology, at plant <- electricity
$ws2.Cells.Item(7,1).Value = "newdb"
$ws2.Cells.Item(7,2).Value = "fancy chemical, from a technology, at plant"
$ws2.Cells.Item(7,3).Value = "newdb"
$ws2.Cells.Item(7,4).Value = "electricity"
$ws2.Cells.Item(7,5).Value = 3
$ws2.Cells.Item(7,6).Value = "kWh"
$ws2.Cells.Item(7,7).Value = "US"

# Row 8: fancy chemical, from a technology, at plant <- 2-Methyl pentane, from chemistry
$ws2.Cells.Item(8,1).Value = "newdb"
$ws2.Cells.Item(8,2).Value = "fancy chemical, from a technology, at plant"
$ws2.Cells.Item(8,3).Value = "newdb"
$ws2.Cells.Item(8,4).Value = "2-Methyl pentane, from chemistry"
$ws2.Cells.Item(8,5).Value = 0.78
$ws2.Cells.Item(8,6).Value = "kg"
$ws2.Cells.Item(8,7).Value = "US"

# Best-effort column widths
$ws2.Columns.Item(1).ColumnWidth = 9
$ws2.Columns.Item(2).ColumnWidth = 39.71
$ws2.Columns.Item(3).ColumnWidth = 12.71
$ws2.Columns.Item(4).ColumnWidth = 32
$ws2.Columns.Item(5).ColumnWidth = 7.86
$ws2.Columns.Item(6).ColumnWidth = 5
$ws2.Columns.Item(7).Delete()
$ws2.Columns.Item(10).Delete()

# ---------------------------------------------------------------------------
# Sheet 3: "Delete Exchanges" -- data unchanged, selection only
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Delete Exchanges")

# ---------------------------------------------------------------------------
# Selections / active sheet
# ---------------------------------------------------------------------------
$ws3.Range("A2").Select()
$ws2.Range("D13").Select()
$ws1.Activate()
$ws1.Range("D6").Select()
